# PM04 Tidsregistrering for Toke.xlsx - add a new day's time entries (2020-05-20)
# to the "Ark1" time-log sheet: one entry continuing the existing table at row 32,
# plus four more entries in rows 33-36, then re-apply the sheet's row-height /
# date-validation formatting to the rows pushed down below the new entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 32: continues the existing log (same row/style as the template used
#     to auto-sum G/H already carried a formula here) ---
$ws.Cells.Item(32, 1).Value = "SD0102+DCD0102"
$ws.Cells.Item(32, 3).Value = 43971
$ws.Cells.Item(32, 4).Value = 0.354166666666667
$ws.Cells.Item(32, 5).Value = 0.416666666666667

# --- Rows 33-36: four more tasks logged the same day ---
$ws.Cells.Item(33, 1).Value = "AD03"
$ws.Cells.Item(33, 3).Value = 43971
$ws.Cells.Item(33, 4).Value = 0.416666666666667
$ws.Cells.Item(33, 5).Value = 0.541666666666667

$ws.Cells.Item(34, 1).Value = "Vejledning fra Karsten og Andárs"
$ws.Cells.Item(34, 3).Value = 43971
$ws.Cells.Item(34, 4).Value = 0.541666666666667
$ws.Cells.Item(34, 5).Value = 0.604166666666667

$ws.Cells.Item(35, 1).Value = "DOM10"
$ws.Cells.Item(35, 3).Value = 43971
$ws.Cells.Item(35, 4).Value = 0.604166666666667
$ws.Cells.Item(35, 5).Value = 0.635416666666667

$ws.Cells.Item(36, 1).Value = "Vejledning fra Anders"
$ws.Cells.Item(36, 3).Value = 43971
$ws.Cells.Item(36, 4).Value = 0.635416666666667
$ws.Cells.Item(36, 5).Value = 0.708333333333333

# --- Row heights: row 32 keeps the "filled" row height, the newly-typed rows
#     33-36 plus all the still-blank template rows pushed further down
#     (37-143) pick up the slightly shorter default height ---
$ws.Range("A32").RowHeight = 18.55
$ws.Range("A33:A143").RowHeight = 13.8

# --- Date validation on column C used to cover C32:C1032 in one block; now
#     that rows 32 (and 33-143) carry real/blank entries, split that block so
#     the still-open date range starts at C144, and give the now-used C32 /
#     C33:C143 ranges their own (wider) date bounds ---
$rngC32 = $ws.Range("C32")
$rngC32.Validation.Delete()
$rngC32.Validation.Add(4, 1, 1, "43881", "1E+19")

$rngC33to143 = $ws.Range("C33:C143")
$rngC33to143.Validation.Delete()
$rngC33to143.Validation.Add(4, 1, 1, "43881", "439080")

# --- Selection left where the user was last working ---
$ws.Range("C22").Select()

Write-Output "done"
